# Reorder the "Recorded By" (column G) comma-separated author lists so
# that any "System"/"system" entries are moved to the end of the list,
# while preserving the relative order of all other entries.
#
# e.g. "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com"     -> "backup@backdoor.com, System, system"
#
# Rows whose value has no "System"/"system" token, or where the
# non-system token(s) already precede it, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $others = @()
    $systems = @()
    foreach ($p in $parts) {
        if ($p -eq "System" -or $p -eq "system") {
            $systems += $p
        } else {
            $others += $p
        }
    }

    if ($systems.Count -eq 0 -or $others.Count -eq 0) {
        continue
    }

    $newParts = $others + $systems
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
